$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()
$newRow.Range.Item(1,1).Value = "SWRS_BRMTR_006"
$newRow.Range.Item(1,2).Value = "Images shall be displayed for 500ms. After this they should disappear, and a grey backround shall be visible."
$newRow.Range.Item(1,3).Value = "Requirement"
$newRow.Range.Item(1,4).Value = "Open"
